$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "windspeed"
$ws.Range("B1").Value = "kW/h"

$ws.Range("E7").Select()
